$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dmstc Stndrd Upld Tmplt")

# --- Row 5: AC5 loses its "wei:0" label value, just keeping the row's
#     background/border formatting (matches the neighboring AD5 cell).
$ws.Range("AC5").ClearContents()
$ws.Range("AD5").Copy()
$ws.Range("AC5").PasteSpecial(-4122)  # xlPasteFormats

# --- Rows 6-23: replace the old sample numeric/text filler data with the
#     new "bean generator" test values (vpn/ven/lab/siz/typ/pon/cla/cat/mat).
#     First wipe the old filler values (B:M) for rows 6-41 (the max extent
#     of the old sample data) so no stale values/cells linger in
#     columns/rows that no longer get new data.
$ws.Range("B6:M41").Clear()

$vpn = @("vpn:0","vpn:1","vpn:2","vpn:3","vpn:4","vpn:5","vpn:6","vpn:7","vpn:8","vpn:9","vpn:10")
$ven = @("ven:0","ven:1","ven:2","ven:3","ven:4","ven:5","ven:6","ven:7","ven:8","ven:9","ven:10","ven:11","ven:12","ven:13")
$lab = @("lab:0","lab:1","lab:2","lab:3","lab:4","lab:5","lab:6","lab:7","lab:8")
$siz = @("siz:0","siz:1","siz:2","siz:3","siz:4","siz:5","siz:6","siz:7","siz:8","siz:9","siz:10","siz:11","siz:12","siz:13","siz:14","siz:15","siz:16")
$typ = @("typ:0","typ:1","typ:2","typ:3","typ:4")
$pon = @("pon:0","pon:1","pon:2","pon:3","pon:4","pon:5","pon:6","pon:7","pon:8","pon:9","pon:10")
$cla = @("cla:0","cla:1","cla:2","cla:3")
$cat = @("cat:0","cat:1","cat:2","cat:3","cat:4","cat:5","cat:6","cat:7","cat:8")
$mat = @("mat:0","mat:1","mat:2","mat:3")

$colE = @(43.328201788639205,34.044522735640385,21.983320925782913,92.67395672507635,71.11132427952623,46.64655779529334,33.04470481434443,34.80922416483404,46.23856997593322,36.73128468991719,20.003075783816758,13.095175229733357,36.87265333081227,18.23226929096021,27.121965955706816,21.940724113431877,14.723589109087342,84.06986321340636)
$colI = @(30.882660345501446,99.5077370191724,75.84210983404257,0.1640512395721161,95.61802755902964)
$colJ = @(64.42700747642085,4.109921346892486,82.8355499140805,37.16781275185029,74.35902886637639,6.371942900982064,86.27026142499349,97.90472142509395,77.92873647610753,90.96274319716811,1.47455171339691)

for ($i = 0; $i -lt 18; $i++) {
    $row = 6 + $i

    if ($i -lt $vpn.Count) { $ws.Cells.Item($row, 2).Value2 = $vpn[$i] }   # B
    if ($i -lt $ven.Count) { $ws.Cells.Item($row, 3).Value2 = $ven[$i] }   # C
    if ($i -lt $lab.Count) { $ws.Cells.Item($row, 4).Value2 = $lab[$i] }   # D

    $ws.Cells.Item($row, 5).Value2 = $colE[$i]                             # E
    $ws.Cells.Item($row, 6).Value2 = $siz[$i]                              # F

    if ($i -lt $typ.Count) { $ws.Cells.Item($row, 7).Value2 = $typ[$i] }   # G
    if ($i -lt $pon.Count) { $ws.Cells.Item($row, 8).Value2 = $pon[$i] }   # H
    if ($i -lt $colI.Count) { $ws.Cells.Item($row, 9).Value2 = $colI[$i] } # I

    $ws.Cells.Item($row, 10).Value2 = $colJ[$i]                            # J

    if ($i -lt $cla.Count) { $ws.Cells.Item($row, 11).Value2 = $cla[$i] }  # K
    if ($i -lt $cat.Count) { $ws.Cells.Item($row, 12).Value2 = $cat[$i] }  # L
    if ($i -lt $mat.Count) { $ws.Cells.Item($row, 13).Value2 = $mat[$i] }  # M
}

# --- Clear the old AC6:AC41 "wei/sup" helper-list labels entirely (column
#     no longer used as a lookup source for these rows).
$ws.Range("AC6:AC46").Clear()

# --- Rows 24-46 no longer exist; the template now only spans to row 23.
$ws.Rows("24:46").Delete()
